$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: 15.5546875 -> 14.5546875 (closest achievable via ColumnWidth quantization is 14.5)
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666

# Write updated data values (A,B,C) for rows 1-250
$arr = New-Object 'object[,]' 250,3
$arr[0,0]=0.0; $arr[0,1]=0.043179999999999996; $arr[0,2]=6.612813525522457
$arr[1,0]=0.0009588755020080321; $arr[1,1]=0.043179999999999996; $arr[1,2]=6.612813525522457
$arr[2,0]=0.0019177510040160643; $arr[2,1]=0.043179999999999996; $arr[2,2]=6.612813525522457
$arr[3,0]=0.0028766265060240965; $arr[3,1]=0.043179999999999996; $arr[3,2]=6.612813525522457
$arr[4,0]=0.0038355020080321285; $arr[4,1]=0.043179999999999996; $arr[4,2]=6.612813525522457
$arr[5,0]=0.004794377510040161; $arr[5,1]=0.043179999999999996; $arr[5,2]=6.612813525522457
$arr[6,0]=0.005753253012048193; $arr[6,1]=0.043179999999999996; $arr[6,2]=6.612813525522457
$arr[7,0]=0.006712128514056225; $arr[7,1]=0.043179999999999996; $arr[7,2]=6.612813525522457
$arr[8,0]=0.007671004016064257; $arr[8,1]=0.043179999999999996; $arr[8,2]=6.612813525522457
$arr[9,0]=0.008629879518072288; $arr[9,1]=0.043179999999999996; $arr[9,2]=6.612813525522457
$arr[10,0]=0.009588755020080322; $arr[10,1]=0.043179999999999996; $arr[10,2]=6.612813525522457
$arr[11,0]=0.010547630522088352; $arr[11,1]=0.043179999999999996; $arr[11,2]=6.612813525522457
$arr[12,0]=0.011506506024096386; $arr[12,1]=0.043179999999999996; $arr[12,2]=6.612813525522457
$arr[13,0]=0.012465381526104416; $arr[13,1]=0.043179999999999996; $arr[13,2]=6.612813525522457
$arr[14,0]=0.01342425702811245; $arr[14,1]=0.043179999999999996; $arr[14,2]=6.612813525522457
$arr[15,0]=0.01438313253012048; $arr[15,1]=0.043179999999999996; $arr[15,2]=6.612813525522457
$arr[16,0]=0.015342008032128514; $arr[16,1]=0.043179999999999996; $arr[16,2]=6.612813525522457
$arr[17,0]=0.016300883534136546; $arr[17,1]=0.043179999999999996; $arr[17,2]=6.612813525522457
$arr[18,0]=0.017259759036144576; $arr[18,1]=0.043179999999999996; $arr[18,2]=6.612813525522457
$arr[19,0]=0.01821863453815261; $arr[19,1]=0.043179999999999996; $arr[19,2]=6.612813525522457
$arr[20,0]=0.019177510040160644; $arr[20,1]=0.043179999999999996; $arr[20,2]=6.612813525522457
$arr[21,0]=0.020136385542168674; $arr[21,1]=0.043179999999999996; $arr[21,2]=6.612813525522457
$arr[22,0]=0.021095261044176705; $arr[22,1]=0.043179999999999996; $arr[22,2]=6.612813525522457
$arr[23,0]=0.02205413654618474; $arr[23,1]=0.043179999999999996; $arr[23,2]=6.612813525522457
$arr[24,0]=0.023013012048192772; $arr[24,1]=0.043179999999999996; $arr[24,2]=6.612813525522457
$arr[25,0]=0.023971887550200806; $arr[25,1]=0.043179999999999996; $arr[25,2]=6.612813525522457
$arr[26,0]=0.024930763052208833; $arr[26,1]=0.043179999999999996; $arr[26,2]=6.612813525522457
$arr[27,0]=0.025889638554216866; $arr[27,1]=0.043179999999999996; $arr[27,2]=6.612813525522457
$arr[28,0]=0.0268485140562249; $arr[28,1]=0.043179999999999996; $arr[28,2]=6.612813525522457
$arr[29,0]=0.02780738955823293; $arr[29,1]=0.043179999999999996; $arr[29,2]=6.612813525522457
$arr[30,0]=0.02876626506024096; $arr[30,1]=0.043179999999999996; $arr[30,2]=6.612813525522457
$arr[31,0]=0.029725140562248994; $arr[31,1]=0.043179999999999996; $arr[31,2]=6.612813525522457
$arr[32,0]=0.030684016064257028; $arr[32,1]=0.043179999999999996; $arr[32,2]=6.612813525522457
$arr[33,0]=0.031642891566265055; $arr[33,1]=0.043179999999999996; $arr[33,2]=6.612813525522457
$arr[34,0]=0.03260176706827309; $arr[34,1]=0.043179999999999996; $arr[34,2]=6.612813525522457
$arr[35,0]=0.03356064257028112; $arr[35,1]=0.043179999999999996; $arr[35,2]=6.612813525522457
$arr[36,0]=0.03451951807228915; $arr[36,1]=0.043179999999999996; $arr[36,2]=6.612813525522457
$arr[37,0]=0.03547839357429718; $arr[37,1]=0.043179999999999996; $arr[37,2]=6.612813525522457
$arr[38,0]=0.03643726907630522; $arr[38,1]=0.043179999999999996; $arr[38,2]=6.612813525522457
$arr[39,0]=0.03739614457831325; $arr[39,1]=0.043179999999999996; $arr[39,2]=6.612813525522457
$arr[40,0]=0.03835502008032129; $arr[40,1]=0.043179999999999996; $arr[40,2]=6.612813525522457
$arr[41,0]=0.03931389558232932; $arr[41,1]=0.043179999999999996; $arr[41,2]=6.612813525522457
$arr[42,0]=0.04027277108433735; $arr[42,1]=0.043179999999999996; $arr[42,2]=6.612813525522457
$arr[43,0]=0.04123164658634538; $arr[43,1]=0.043179999999999996; $arr[43,2]=6.612813525522457
$arr[44,0]=0.04219052208835341; $arr[44,1]=0.043179999999999996; $arr[44,2]=6.612813525522457
$arr[45,0]=0.043149397590361446; $arr[45,1]=0.043179999999999996; $arr[45,2]=6.612813525522457
$arr[46,0]=0.04410827309236948; $arr[46,1]=0.043179999999999996; $arr[46,2]=6.612813525522457
$arr[47,0]=0.04506714859437751; $arr[47,1]=0.043179999999999996; $arr[47,2]=6.612813525522457
$arr[48,0]=0.046026024096385544; $arr[48,1]=0.043179999999999996; $arr[48,2]=6.612813525522457
$arr[49,0]=0.046984899598393574; $arr[49,1]=0.043179999999999996; $arr[49,2]=6.612813525522457
$arr[50,0]=0.04794377510040161; $arr[50,1]=0.043179999999999996; $arr[50,2]=6.612813525522457
$arr[51,0]=0.048902650602409635; $arr[51,1]=0.043179999999999996; $arr[51,2]=6.612813525522457
$arr[52,0]=0.049861526104417665; $arr[52,1]=0.043179999999999996; $arr[52,2]=6.612813525522457
$arr[53,0]=0.0508204016064257; $arr[53,1]=0.043179999999999996; $arr[53,2]=6.612813525522457
$arr[54,0]=0.05177927710843373; $arr[54,1]=0.043179999999999996; $arr[54,2]=6.612813525522457
$arr[55,0]=0.05273815261044176; $arr[55,1]=0.043179999999999996; $arr[55,2]=6.612813525522457
$arr[56,0]=0.0536970281124498; $arr[56,1]=0.043179999999999996; $arr[56,2]=6.612813525522457
$arr[57,0]=0.05465590361445783; $arr[57,1]=0.043179999999999996; $arr[57,2]=6.612813525522457
$arr[58,0]=0.05561477911646586; $arr[58,1]=0.043179999999999996; $arr[58,2]=6.612813525522457
$arr[59,0]=0.05657365461847389; $arr[59,1]=0.043179999999999996; $arr[59,2]=6.612813525522457
$arr[60,0]=0.05753253012048192; $arr[60,1]=0.043179999999999996; $arr[60,2]=6.612813525522457
$arr[61,0]=0.05849140562248996; $arr[61,1]=0.043179999999999996; $arr[61,2]=6.612813525522457
$arr[62,0]=0.05945028112449799; $arr[62,1]=0.043179999999999996; $arr[62,2]=6.612813525522457
$arr[63,0]=0.060409156626506026; $arr[63,1]=0.043179999999999996; $arr[63,2]=6.612813525522457
$arr[64,0]=0.061368032128514056; $arr[64,1]=0.043179999999999996; $arr[64,2]=6.612813525522457
$arr[65,0]=0.06232690763052209; $arr[65,1]=0.043179999999999996; $arr[65,2]=6.612813525522457
$arr[66,0]=0.06328578313253011; $arr[66,1]=0.043179999999999996; $arr[66,2]=6.612813525522457
$arr[67,0]=0.06424465863453814; $arr[67,1]=0.043179999999999996; $arr[67,2]=6.612813525522457
$arr[68,0]=0.06520353413654618; $arr[68,1]=0.043179999999999996; $arr[68,2]=6.612813525522457
$arr[69,0]=0.06616240963855421; $arr[69,1]=0.043179999999999996; $arr[69,2]=6.612813525522457
$arr[70,0]=0.06712128514056225; $arr[70,1]=0.043179999999999996; $arr[70,2]=6.612813525522457
$arr[71,0]=0.06808016064257028; $arr[71,1]=0.043179999999999996; $arr[71,2]=6.612813525522457
$arr[72,0]=0.0690390361445783; $arr[72,1]=0.043179999999999996; $arr[72,2]=6.612813525522457
$arr[73,0]=0.06999791164658634; $arr[73,1]=0.043179999999999996; $arr[73,2]=6.612813525522457
$arr[74,0]=0.07095678714859437; $arr[74,1]=0.043179999999999996; $arr[74,2]=6.612813525522457
$arr[75,0]=0.0719156626506024; $arr[75,1]=0.043179999999999996; $arr[75,2]=6.612813525522457
$arr[76,0]=0.07287453815261044; $arr[76,1]=0.043179999999999996; $arr[76,2]=6.612813525522457
$arr[77,0]=0.07383341365461847; $arr[77,1]=0.043179999999999996; $arr[77,2]=6.612813525522457
$arr[78,0]=0.0747922891566265; $arr[78,1]=0.043179999999999996; $arr[78,2]=6.612813525522457
$arr[79,0]=0.07575116465863453; $arr[79,1]=0.043179999999999996; $arr[79,2]=6.612813525522457
$arr[80,0]=0.07671004016064258; $arr[80,1]=0.043179999999999996; $arr[80,2]=6.612813525522457
$arr[81,0]=0.07766891566265059; $arr[81,1]=0.043179999999999996; $arr[81,2]=6.612813525522457
$arr[82,0]=0.07862779116465864; $arr[82,1]=0.043179999999999996; $arr[82,2]=6.612813525522457
$arr[83,0]=0.07958666666666667; $arr[83,1]=0.043179999999999996; $arr[83,2]=6.612813525522457
$arr[84,0]=0.0805455421686747; $arr[84,1]=0.043179999999999996; $arr[84,2]=6.612813525522457
$arr[85,0]=0.08150441767068271; $arr[85,1]=0.043179999999999996; $arr[85,2]=6.612813525522457
$arr[86,0]=0.08246329317269076; $arr[86,1]=0.043179999999999996; $arr[86,2]=6.612813525522457
$arr[87,0]=0.08342216867469877; $arr[87,1]=0.043179999999999996; $arr[87,2]=6.612813525522457
$arr[88,0]=0.08438104417670682; $arr[88,1]=0.043179999999999996; $arr[88,2]=6.612813525522457
$arr[89,0]=0.08533991967871485; $arr[89,1]=0.043179999999999996; $arr[89,2]=6.612813525522457
$arr[90,0]=0.08629879518072289; $arr[90,1]=0.043179999999999996; $arr[90,2]=6.612813525522457
$arr[91,0]=0.08725767068273091; $arr[91,1]=0.043179999999999996; $arr[91,2]=6.612813525522457
$arr[92,0]=0.08821654618473895; $arr[92,1]=0.043179999999999996; $arr[92,2]=6.612813525522457
$arr[93,0]=0.08917542168674698; $arr[93,1]=0.043179999999999996; $arr[93,2]=6.612813525522457
$arr[94,0]=0.09013429718875501; $arr[94,1]=0.043179999999999996; $arr[94,2]=6.612813525522457
$arr[95,0]=0.09109317269076304; $arr[95,1]=0.043179999999999996; $arr[95,2]=6.612813525522457
$arr[96,0]=0.09205204819277109; $arr[96,1]=0.043179999999999996; $arr[96,2]=6.612813525522457
$arr[97,0]=0.0930109236947791; $arr[97,1]=0.043179999999999996; $arr[97,2]=6.612813525522457
$arr[98,0]=0.09396979919678715; $arr[98,1]=0.043179999999999996; $arr[98,2]=6.612813525522457
$arr[99,0]=0.09492867469879518; $arr[99,1]=0.043179999999999996; $arr[99,2]=6.612813525522457
$arr[100,0]=0.09588755020080322; $arr[100,1]=0.043179999999999996; $arr[100,2]=6.612813525522457
$arr[101,0]=0.09684642570281123; $arr[101,1]=0.043179999999999996; $arr[101,2]=6.612813525522457
$arr[102,0]=0.09780530120481927; $arr[102,1]=0.043179999999999996; $arr[102,2]=6.612813525522457
$arr[103,0]=0.0987641767068273; $arr[103,1]=0.043179999999999996; $arr[103,2]=6.612813525522457
$arr[104,0]=0.09972305220883533; $arr[104,1]=0.043179999999999996; $arr[104,2]=6.612813525522457
$arr[105,0]=0.10068192771084336; $arr[105,1]=0.043179999999999996; $arr[105,2]=6.612813525522457
$arr[106,0]=0.1016408032128514; $arr[106,1]=0.043179999999999996; $arr[106,2]=6.612813525522457
$arr[107,0]=0.10259967871485942; $arr[107,1]=0.043179999999999996; $arr[107,2]=6.612813525522457
$arr[108,0]=0.10355855421686747; $arr[108,1]=0.043179999999999996; $arr[108,2]=6.612813525522457
$arr[109,0]=0.1045174297188755; $arr[109,1]=0.043179999999999996; $arr[109,2]=6.612813525522457
$arr[110,0]=0.10547630522088353; $arr[110,1]=0.043179999999999996; $arr[110,2]=6.612813525522457
$arr[111,0]=0.10643518072289156; $arr[111,1]=0.043179999999999996; $arr[111,2]=6.612813525522457
$arr[112,0]=0.1073940562248996; $arr[112,1]=0.043179999999999996; $arr[112,2]=6.612813525522457
$arr[113,0]=0.10835293172690763; $arr[113,1]=0.043179999999999996; $arr[113,2]=6.612813525522457
$arr[114,0]=0.10931180722891566; $arr[114,1]=0.043179999999999996; $arr[114,2]=6.612813525522457
$arr[115,0]=0.11027068273092369; $arr[115,1]=0.043179999999999996; $arr[115,2]=6.612813525522457
$arr[116,0]=0.11122955823293172; $arr[116,1]=0.043179999999999996; $arr[116,2]=6.612813525522457
$arr[117,0]=0.11218843373493974; $arr[117,1]=0.043179999999999996; $arr[117,2]=6.612813525522457
$arr[118,0]=0.11314730923694778; $arr[118,1]=0.043179999999999996; $arr[118,2]=6.612813525522457
$arr[119,0]=0.11410618473895581; $arr[119,1]=0.043179999999999996; $arr[119,2]=6.612813525522457
$arr[120,0]=0.11506506024096384; $arr[120,1]=0.043179999999999996; $arr[120,2]=6.612813525522457
$arr[121,0]=0.11602393574297187; $arr[121,1]=0.043179999999999996; $arr[121,2]=6.612813525522457
$arr[122,0]=0.11698281124497992; $arr[122,1]=0.043179999999999996; $arr[122,2]=6.612813525522457
$arr[123,0]=0.11794168674698795; $arr[123,1]=0.043179999999999996; $arr[123,2]=6.612813525522457
$arr[124,0]=0.11890056224899598; $arr[124,1]=0.043179999999999996; $arr[124,2]=6.612813525522457
$arr[125,0]=0.11985943775100401; $arr[125,1]=0.043179999999999996; $arr[125,2]=6.612813525522457
$arr[126,0]=0.12081831325301205; $arr[126,1]=0.043179999999999996; $arr[126,2]=6.612813525522457
$arr[127,0]=0.12177718875502007; $arr[127,1]=0.043179999999999996; $arr[127,2]=6.612813525522457
$arr[128,0]=0.12273606425702811; $arr[128,1]=0.043179999999999996; $arr[128,2]=6.612813525522457
$arr[129,0]=0.12369493975903614; $arr[129,1]=0.043179999999999996; $arr[129,2]=6.612813525522457
$arr[130,0]=0.12465381526104417; $arr[130,1]=0.043179999999999996; $arr[130,2]=6.612813525522457
$arr[131,0]=0.1256126907630522; $arr[131,1]=0.043179999999999996; $arr[131,2]=6.612813525522457
$arr[132,0]=0.12657156626506022; $arr[132,1]=0.043179999999999996; $arr[132,2]=6.612813525522457
$arr[133,0]=0.12753044176706826; $arr[133,1]=0.043179999999999996; $arr[133,2]=6.612813525522457
$arr[134,0]=0.12848931726907628; $arr[134,1]=0.043179999999999996; $arr[134,2]=6.612813525522457
$arr[135,0]=0.12944819277108432; $arr[135,1]=0.043179999999999996; $arr[135,2]=6.612813525522457
$arr[136,0]=0.13040706827309237; $arr[136,1]=0.043179999999999996; $arr[136,2]=6.612813525522457
$arr[137,0]=0.13136594377510039; $arr[137,1]=0.043179999999999996; $arr[137,2]=6.612813525522457
$arr[138,0]=0.13232481927710843; $arr[138,1]=0.043179999999999996; $arr[138,2]=6.612813525522457
$arr[139,0]=0.13328369477911647; $arr[139,1]=0.043179999999999996; $arr[139,2]=6.612813525522457
$arr[140,0]=0.1342425702811245; $arr[140,1]=0.043179999999999996; $arr[140,2]=6.612813525522457
$arr[141,0]=0.13520144578313253; $arr[141,1]=0.043179999999999996; $arr[141,2]=6.612813525522457
$arr[142,0]=0.13616032128514055; $arr[142,1]=0.043179999999999996; $arr[142,2]=6.612813525522457
$arr[143,0]=0.13711919678714857; $arr[143,1]=0.043179999999999996; $arr[143,2]=6.612813525522457
$arr[144,0]=0.1380780722891566; $arr[144,1]=0.043179999999999996; $arr[144,2]=6.612813525522457
$arr[145,0]=0.13903694779116466; $arr[145,1]=0.043179999999999996; $arr[145,2]=6.612813525522457
$arr[146,0]=0.13999582329317267; $arr[146,1]=0.043179999999999996; $arr[146,2]=6.612813525522457
$arr[147,0]=0.14095469879518072; $arr[147,1]=0.043179999999999996; $arr[147,2]=6.612813525522457
$arr[148,0]=0.14191357429718873; $arr[148,1]=0.043179999999999996; $arr[148,2]=6.612813525522457
$arr[149,0]=0.14287244979919675; $arr[149,1]=0.043179999999999996; $arr[149,2]=6.612813525522457
$arr[150,0]=0.1438313253012048; $arr[150,1]=0.043179999999999996; $arr[150,2]=6.612813525522457
$arr[151,0]=0.14479020080321284; $arr[151,1]=0.043179999999999996; $arr[151,2]=6.612813525522457
$arr[152,0]=0.14574907630522088; $arr[152,1]=0.043179999999999996; $arr[152,2]=6.612813525522457
$arr[153,0]=0.14670795180722893; $arr[153,1]=0.043179999999999996; $arr[153,2]=6.612813525522457
$arr[154,0]=0.14766682730923694; $arr[154,1]=0.043179999999999996; $arr[154,2]=6.612813525522457
$arr[155,0]=0.14862570281124496; $arr[155,1]=0.043179999999999996; $arr[155,2]=6.612813525522457
$arr[156,0]=0.149584578313253; $arr[156,1]=0.043179999999999996; $arr[156,2]=6.612813525522457
$arr[157,0]=0.15054345381526105; $arr[157,1]=0.043179999999999996; $arr[157,2]=6.612813525522457
$arr[158,0]=0.15150232931726906; $arr[158,1]=0.04316639456547871; $arr[158,2]=6.608646965669487
$arr[159,0]=0.1524612048192771; $arr[159,1]=0.0430854789636675; $arr[159,2]=6.58389430967362
$arr[160,0]=0.15342008032128515; $arr[160,1]=0.042930409964303105; $arr[160,2]=6.536587390422633
$arr[161,0]=0.15437895582329314; $arr[161,1]=0.04269836454679061; $arr[161,2]=6.466115856432249
$arr[162,0]=0.15533783132530118; $arr[162,1]=0.04238484236185508; $arr[162,2]=6.371506709182842
$arr[163,0]=0.15629670682730923; $arr[163,1]=0.041983183327508354; $arr[163,2]=6.251320004901346
$arr[164,0]=0.15725558232931727; $arr[164,1]=0.04148374014665972; $arr[164,2]=6.103469948638477
$arr[165,0]=0.15821445783132532; $arr[165,1]=0.040872424322912; $arr[165,2]=5.9249105354508975
$arr[166,0]=0.15917333333333333; $arr[166,1]=0.04012799382996244; $arr[166,2]=5.711049130659882
$arr[167,0]=0.16013220883534138; $arr[167,1]=0.0392229313296554; $arr[167,2]=5.45633586387562
$arr[168,0]=0.1610910843373494; $arr[168,1]=0.038264055827647386; $arr[168,2]=5.192816824679169
$arr[169,0]=0.16204995983935744; $arr[169,1]=0.03730518032563934; $arr[169,2]=4.935819705104781
$arr[170,0]=0.16300883534136543; $arr[170,1]=0.036346304823631326; $arr[170,2]=4.685344505152472
$arr[171,0]=0.16396771084337347; $arr[171,1]=0.03538742932162331; $arr[171,2]=4.441391224822232
$arr[172,0]=0.16492658634538152; $arr[172,1]=0.034428553819615265; $arr[172,2]=4.203959864114058
$arr[173,0]=0.16588546184738953; $arr[173,1]=0.03346967831760722; $arr[173,2]=3.9730504230279524
$arr[174,0]=0.16684433734939755; $arr[174,1]=0.03251080281559918; $arr[174,2]=3.748662901563918
$arr[175,0]=0.1678032128514056; $arr[175,1]=0.03155192731359116; $arr[175,2]=3.530797299721962
$arr[176,0]=0.16876208835341364; $arr[176,1]=0.030593051811583116; $arr[176,2]=3.319453617502069
$arr[177,0]=0.16972096385542165; $arr[177,1]=0.0296341763095751; $arr[177,2]=3.1146318549042533
$arr[178,0]=0.1706798393574297; $arr[178,1]=0.028675300807567056; $arr[178,2]=2.9163320119285023
$arr[179,0]=0.17163871485943774; $arr[179,1]=0.02771642530555904; $arr[179,2]=2.724554088574828
$arr[180,0]=0.17259759036144579; $arr[180,1]=0.026757549803550995; $arr[180,2]=2.5392980848432183
$arr[181,0]=0.17355646586345377; $arr[181,1]=0.02579867430154295; $arr[181,2]=2.3605640007336794
$arr[182,0]=0.17451534136546182; $arr[182,1]=0.024839798799534935; $arr[182,2]=2.188351836246217
$arr[183,0]=0.17547421686746986; $arr[183,1]=0.02388092329752689; $arr[183,2]=2.0226615913808192
$arr[184,0]=0.1764330923694779; $arr[184,1]=0.022922047795518874; $arr[184,2]=1.8634932661374972
$arr[185,0]=0.17739196787148592; $arr[185,1]=0.02196317229351083; $arr[185,2]=1.7108468605162417
$arr[186,0]=0.17835084337349397; $arr[186,1]=0.021004296791502813; $arr[186,2]=1.5647223745170609
$arr[187,0]=0.179309718875502; $arr[187,1]=0.02004542128949477; $arr[187,2]=1.4251198081399468
$arr[188,0]=0.18026859437751003; $arr[188,1]=0.019145340378242136; $arr[188,2]=1.3000114698923837
$arr[189,0]=0.18122746987951804; $arr[189,1]=0.018437059540459856; $arr[189,2]=1.2056029906054122
$arr[190,0]=0.1821863453815261; $arr[190,1]=0.01788421113480837; $arr[190,2]=1.134385265770936
$arr[191,0]=0.18314522088353413; $arr[191,1]=0.0174595910573836; $arr[191,2]=1.0811579220710357
$arr[192,0]=0.18410409638554218; $arr[192,1]=0.017146522161506472; $arr[192,2]=1.042732933577169
$arr[193,0]=0.18506297188755022; $arr[193,1]=0.01693453756416042; $arr[193,2]=1.0171094367669415
$arr[194,0]=0.1860218473895582; $arr[194,1]=0.016817290602423818; $arr[194,2]=1.0030741949672455
$arr[195,0]=0.18698072289156625; $arr[195,1]=0.016791500235609227; $arr[195,2]=1.0
$arr[196,0]=0.1879395983935743; $arr[196,1]=0.01685646534521964; $arr[196,2]=1.0077528251044863
$arr[197,0]=0.18889847389558231; $arr[197,1]=0.01701396604816778; $arr[197,2]=1.0266729598387934
$arr[198,0]=0.18985734939759036; $arr[198,1]=0.01724091232764501; $arr[198,2]=1.0542448473181347
$arr[199,0]=0.1908162248995984; $arr[199,1]=0.01747111796827942; $arr[199,2]=1.0825859652439427
$arr[200,0]=0.19177510040160645; $arr[200,1]=0.017701323608913837; $arr[200,2]=1.1113029923374258
$arr[201,0]=0.1927339759036144; $arr[201,1]=0.01793152924954824; $arr[201,2]=1.140395928598581
$arr[202,0]=0.19369285140562245; $arr[202,1]=0.01816173489018266; $arr[202,2]=1.1698647740274113
$arr[203,0]=0.1946517269076305; $arr[203,1]=0.018391940530817063; $arr[203,2]=1.1997095286239143
$arr[204,0]=0.19561060240963854; $arr[204,1]=0.018622146171451474; $arr[204,2]=1.229930192388092
$arr[205,0]=0.19656947791164656; $arr[205,1]=0.01885235181208589; $arr[205,2]=1.2605267653199448
$arr[206,0]=0.1975283534136546; $arr[206,1]=0.019082557452720295; $arr[206,2]=1.2914992474194693
$arr[207,0]=0.19848722891566262; $arr[207,1]=0.019312763093354713; $arr[207,2]=1.3228476386866694
$arr[208,0]=0.19944610441767066; $arr[208,1]=0.019542968733989124; $arr[208,2]=1.3545719391215427
$arr[209,0]=0.20040497991967868; $arr[209,1]=0.01977317437462354; $arr[209,2]=1.3866721487240912
$arr[210,0]=0.20136385542168672; $arr[210,1]=0.020003380015257945; $arr[210,2]=1.4191482674943112
$arr[211,0]=0.20232273092369477; $arr[211,1]=0.020233585655892363; $arr[211,2]=1.4520002954322069
$arr[212,0]=0.2032816064257028; $arr[212,1]=0.020463791296526767; $arr[212,2]=1.4852282325377748
$arr[213,0]=0.20424048192771083; $arr[213,1]=0.020693996937161178; $arr[213,2]=1.518832078811018
$arr[214,0]=0.20519935742971884; $arr[214,1]=0.020924202577795595; $arr[214,2]=1.5528118342519357
$arr[215,0]=0.2061582329317269; $arr[215,1]=0.02115440821843; $arr[215,2]=1.5871674988605253
$arr[216,0]=0.20711710843373493; $arr[216,1]=0.021384613859064417; $arr[216,2]=1.6218990726367906
$arr[217,0]=0.20807598393574295; $arr[217,1]=0.021614819499698828; $arr[217,2]=1.657006555580729
$arr[218,0]=0.209034859437751; $arr[218,1]=0.021845025140333246; $arr[218,2]=1.6924899476923425
$arr[219,0]=0.20999373493975904; $arr[219,1]=0.02207523078096765; $arr[219,2]=1.7283492489716277
$arr[220,0]=0.21095261044176705; $arr[220,1]=0.022305436421602067; $arr[220,2]=1.7645844594185887
$arr[221,0]=0.21191148594377507; $arr[221,1]=0.02253564206223647; $arr[221,2]=1.8011955790332215
$arr[222,0]=0.2128703614457831; $arr[222,1]=0.02276584770287089; $arr[222,2]=1.8381826078155308
$arr[223,0]=0.21382923694779116; $arr[223,1]=0.0229960533435053; $arr[223,2]=1.8755455457655126
$arr[224,0]=0.2147881124497992; $arr[224,1]=0.023226258984139703; $arr[224,2]=1.9132843928831667
$arr[225,0]=0.21574698795180722; $arr[225,1]=0.02345646462477412; $arr[225,2]=1.9513991491684977
$arr[226,0]=0.21670586345381526; $arr[226,1]=0.023686670265408532; $arr[226,2]=1.9898898146215012
$arr[227,0]=0.21766473895582328; $arr[227,1]=0.02391687590604295; $arr[227,2]=2.02875638924218
$arr[228,0]=0.21862361445783132; $arr[228,1]=0.024147081546677353; $arr[228,2]=2.0679988730305294
$arr[229,0]=0.21958248995983934; $arr[229,1]=0.02437728718731177; $arr[229,2]=2.107617265986556
$arr[230,0]=0.22054136546184738; $arr[230,1]=0.024607492827946182; $arr[230,2]=2.1476115681102557
$arr[231,0]=0.22150024096385543; $arr[231,1]=0.024837698468580593; $arr[231,2]=2.1879817794016287
$arr[232,0]=0.22245911646586344; $arr[232,1]=0.025067904109215004; $arr[232,2]=2.228727899860675
$arr[233,0]=0.22341799196787146; $arr[233,1]=0.025298109749849414; $arr[233,2]=2.2698499294873957
$arr[234,0]=0.22437686746987948; $arr[234,1]=0.025528315390483832; $arr[234,2]=2.3113478682817923
$arr[235,0]=0.22533574297188752; $arr[235,1]=0.025758521031118236; $arr[235,2]=2.3532217162438593
$arr[236,0]=0.22629461847389556; $arr[236,1]=0.025988726671752654; $arr[236,2]=2.3954714733736027
$arr[237,0]=0.22725349397590358; $arr[237,1]=0.026218932312387058; $arr[237,2]=2.438097139671018
$arr[238,0]=0.22821236947791163; $arr[238,1]=0.026449137953021475; $arr[238,2]=2.4810987151361097
$arr[239,0]=0.22917124497991967; $arr[239,1]=0.026679343593655886; $arr[239,2]=2.5244761997688743
$arr[240,0]=0.23013012048192769; $arr[240,1]=0.026909549234290297; $arr[240,2]=2.5682295935693116
$arr[241,0]=0.2310889959839357; $arr[241,1]=0.027139754874924708; $arr[241,2]=2.612358896537424
$arr[242,0]=0.23204787148594375; $arr[242,1]=0.02736996051555912; $arr[242,2]=2.65686410867321
$arr[243,0]=0.2330067469879518; $arr[243,1]=0.027600166156193536; $arr[243,2]=2.7017452299766713
$arr[244,0]=0.23396562248995983; $arr[244,1]=0.027830371796827947; $arr[244,2]=2.7470022604478044
$arr[245,0]=0.23492449799196785; $arr[245,1]=0.028060577437462358; $arr[245,2]=2.7926352000866124
$arr[246,0]=0.2358833734939759; $arr[246,1]=0.02829078307809676; $arr[246,2]=2.8386440488930926
$arr[247,0]=0.2368422489959839; $arr[247,1]=0.02852098871873118; $arr[247,2]=2.8850288068672496
$arr[248,0]=0.23780112449799196; $arr[248,1]=0.02875119435936559; $arr[248,2]=2.931789474009079
$arr[249,0]=0.23875999999999997; $arr[249,1]=0.0289814; $arr[249,2]=2.9789260503185817
$ws.Range("A1:C250").Value = $arr
